# "Completed my task for add user." — populate the "Yash" sheet with a new
# add-user test-data table (Employee Name / Username / Password / Confirm
# Password header row + one data row), style it, size the columns, and
# leave that sheet active/selected (mirrors the OrangeHRM add-user test
# data fixture added in the commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yash")

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Employee Name"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Confirm Password"

# --- Data row ---------------------------------------------------------------
$ws.Range("A2").Value = "Jordan"
$ws.Range("B2").Value = "Hodo5613"
$ws.Range("C2").Value = "hodo0219"
$ws.Range("D2").Value = "hodo0219"

# --- Formatting: thin box border around the whole table, yellow fill on the
#     header row -------------------------------------------------------------
$ws.Range("A1:D2").Borders.LineStyle = 1
$ws.Range("A1:D2").Borders.Weight = 2
$ws.Range("A1:D1").Interior.Color = 65535

# --- Column widths (approximate autosize-to-content) -----------------------
$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 20.5

# --- Leave selection on D9 and make this the active sheet/tab --------------
[void]$ws.Range("D9").Select()
[void]$ws.Activate()
